$d = $word.ActiveDocument

# Locate the two paragraphs to remove ("Ver no Jupiter Salvar em pdf
# Salvar em docx" and the "© 2020 ..." footer line) plus the blank
# paragraph that used to separate them from the page-break paragraph,
# by matching on their text rather than hard-coded indices so the
# script is resilient to any paragraphs that may already have shifted.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text

    if ($text -like "Ver no Jupiter*") {
        $startPara = $para
    }
    if ($text -like "*Contact: luizeleno@usp.br*") {
        # The paragraph right after the copyright line is the blank
        # paragraph that should be removed along with it.
        $endPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
